# Apply updated cryptocurrency price/volume data to sheet1 (matches commit diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '90.953.44'
$ws.Range("E2").Value = '  +1.67%  '

$ws.Range("D3").Value = '3.149.77'
$ws.Range("E3").Value = '  +3.10%  '

$ws.Range("E4").Value = '  -0.46%  '

$ws.Range("D5").Value = '''239.70'
$ws.Range("E5").Value = '  +2.10%  '

$ws.Range("D6").Value = '''618.68'
$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("D7").Value = '''1.12'
$ws.Range("E7").Value = '  +6.33%  '

$ws.Range("D8").Value = '''0.374'
$ws.Range("E8").Value = '  +4.15%  '

$ws.Range("E9").Value = '  -0.07%  '

$ws.Range("D10").Value = '3.149.05'
$ws.Range("E10").Value = '  +3.11%  '

$ws.Range("D11").Value = '''0.745'
$ws.Range("E11").Value = '  +6.64%  '

$ws.Range("E12").Value = '  +1.97%  '

$ws.Range("D13").Value = '''0.0000247'
$ws.Range("E13").Value = '  +1.68%  '

$ws.Range("D14").Value = '''35.14'
$ws.Range("E14").Value = '  +2.31%  '

$ws.Range("D15").Value = '''5.57'
$ws.Range("E15").Value = '  +4.54%  '

$ws.Range("D16").Value = '90.921.60'
$ws.Range("E16").Value = '  +1.79%  '

$ws.Range("D17").Value = '3.746.14'

$ws.Range("D18").Value = '3.152.59'
$ws.Range("E18").Value = '  +2.64%  '

$ws.Range("D19").Value = '''3.74'
$ws.Range("E19").Value = '  +0.89%  '

$ws.Range("D20").Value = '''15.17'
$ws.Range("E20").Value = '  +11.21%  '

$ws.Range("D21").Value = '''6.03'
$ws.Range("E21").Value = '  +12.61%  '

$ws.Range("D22").Value = '''450.35'
$ws.Range("E22").Value = '  +5.35%  '

$ws.Range("E23").Value = '  -1.79%  '

$ws.Range("D24").Value = '''9.15'
$ws.Range("E24").Value = '  +6.16%  '

$ws.Range("D25").Value = '''5.74'
$ws.Range("E25").Value = '  +4.21%  '

$ws.Range("D26").Value = '''88.78'
$ws.Range("E26").Value = '  +9.03%  '

$ws.Range("D27").Value = '''12.01'
$ws.Range("E27").Value = '  +4.00%  '

$ws.Range("E28").Value = '  +2.86%  '

$ws.Range("E29").Value = '  -0.18%  '

$ws.Range("D30").Value = '''0.141'
$ws.Range("E30").Value = '  +58.91%  '

$ws.Range("D31").Value = '''0.235'
$ws.Range("E31").Value = '  +19.94%  '

$ws.Range("D32").Value = '''0.171'
$ws.Range("E32").Value = '  +10.83%  '

$ws.Range("D33").Value = '''9.37'
$ws.Range("E33").Value = '  +5.13%  '

$ws.Range("E34").Value = '  +16.09%  '

$ws.Range("E35").Value = '  -5.25%  '

$ws.Range("E36").Value = '  +9.66%  '

$ws.Range("D37").Value = '''26.38'
$ws.Range("E37").Value = '  +4.31%  '

$ws.Range("B38").Value = 'PancakeSwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D38").Value = '''1.97'
$ws.Range("E38").Value = '  +5.79%  '

$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '''508.44'
$ws.Range("E39").Value = '  +4.44%  '

$ws.Range("D40").Value = '''1.34'
$ws.Range("E40").Value = '  +8.04%  '

$ws.Range("D41").Value = '''3.84'
$ws.Range("E41").Value = '  -1.93%  '

$ws.Range("D42").Value = '''0.448'
$ws.Range("E42").Value = '  +13.46%  '

$ws.Range("D43").Value = '''3.45'
$ws.Range("E43").Value = '  -3.36%  '

$ws.Range("D44").Value = '''22.10'
$ws.Range("E44").Value = '  +0.10%  '

$ws.Range("D46").Value = '''0.717'
$ws.Range("E46").Value = '  +7.62%  '

$ws.Range("D47").Value = '''1.94'
$ws.Range("E47").Value = '  +5.63%  '

$ws.Range("D48").Value = '''155.70'
$ws.Range("E48").Value = '  -0.48%  '

$ws.Range("D49").Value = '''1.37'
$ws.Range("E49").Value = '  +7.08%  '

$ws.Range("D50").Value = '''4.47'
$ws.Range("E50").Value = '  +4.38%  '

$ws.Range("B51").Value = 'OKB'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D51").Value = '''44.09'
$ws.Range("E51").Value = '  -0.39%  '
